# "Version for hospitalized patients"
#  - rename header B1 from "ntot" to "nhos"
#  - shift the date series two days earlier: insert 25.02.2020 / 26.02.2020
#    at the top and append 24.03.2020 at the bottom (every other date moves
#    up two rows)
#  - replace the nhos/nicu counts with the new (hospitalized-patients) series
#  - drop the grey shading that used to mark column B and the B13 formula
#  - grow the sheet from 27 to 30 data rows, and move the view back to A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header -----------------------------------------------------------
$ws.Range("B1").Value = "nhos"

# --- the new data, one row per day ------------------------------------
$rows = @(
    @{Row=2;  Date="25.02.2020"; B=4;   C=0},
    @{Row=3;  Date="26.02.2020"; B=4;   C=0},
    @{Row=4;  Date="27.02.2020"; B=4;   C=0},
    @{Row=5;  Date="28.02.2020"; B=4;   C=0},
    @{Row=6;  Date="29.02.2020"; B=4;   C=0},
    @{Row=7;  Date="01.03.2020"; B=4;   C=0},
    @{Row=8;  Date="02.03.2020"; B=6;   C=0},
    @{Row=9;  Date="03.03.2020"; B=8;   C=0},
    @{Row=10; Date="04.03.2020"; B=11;  C=1},
    @{Row=11; Date="05.03.2020"; B=14;  C=1},
    @{Row=12; Date="06.03.2020"; B=15;  C=2},
    @{Row=13; Date="07.03.2020"; B=16;  C=4},
    @{Row=14; Date="08.03.2020"; B=22;  C=3},
    @{Row=15; Date="09.03.2020"; B=29;  C=5},
    @{Row=16; Date="10.03.2020"; B=36;  C=6},
    @{Row=17; Date="11.03.2020"; B=38;  C=7},
    @{Row=18; Date="12.03.2020"; B=43;  C=8},
    @{Row=19; Date="13.03.2020"; B=52;  C=10},
    @{Row=20; Date="14.03.2020"; B=62;  C=14},
    @{Row=21; Date="15.03.2020"; B=78;  C=19},
    @{Row=22; Date="16.03.2020"; B=110; C=27},
    @{Row=23; Date="17.03.2020"; B=117; C=35},
    @{Row=24; Date="18.03.2020"; B=118; C=34},
    @{Row=25; Date="19.03.2020"; B=140; C=32},
    @{Row=26; Date="20.03.2020"; B=152; C=30},
    @{Row=27; Date="21.03.2020"; B=175; C=23},
    @{Row=28; Date="22.03.2020"; B=203; C=23},
    @{Row=29; Date="23.03.2020"; B=223; C=41},
    @{Row=30; Date="24.03.2020"; B=266; C=46}
)

# Put column A into Text format first so the "dd.mm.yyyy" strings go in as
# plain text (matching the rest of the sheet) instead of being auto-parsed
# into date serial numbers; flip it back to General afterwards so the
# cells end up on the same (default) style as before.
$ws.Range("A2:A30").NumberFormat = "@"
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}
$ws.Range("A2:A30").NumberFormat = "General"

# --- drop the grey fill that used to sit on column B (rows 2-13) ------
# Rows 14+ already carry the default style; copy it onto 2-13 so the whole
# column goes back to plain, unshaded cells (the old B13 formula has
# already been overwritten by a literal value above).
$ws.Range("B14").Copy()
$ws.Range("B2:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- view bookkeeping ---------------------------------------------------
$ws.Application.GoTo($ws.Range("A1"), $true)
$ws.Range("C31").Select()
